$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.943.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -1.33%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'3.039.80"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +0.69%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'384.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.39%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'101.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -0.98%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.533"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  -2.54%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.574"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -2.75%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'36.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -1.09%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  -0.16%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.0843"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -2.29%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'3.530.82"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.84%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'18.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -1.62%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'7.62"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -1.28%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.039.04"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +0.71%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'0.974"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -0.01%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'10.62"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.32%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'50.981.82"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -1.25%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'3.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +2.57%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.0₃0948"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -1.63%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'12.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -2.47%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'69.28"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -1.04%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'262.66"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -1.78%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'3.09"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -1.81%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'7.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -5.57%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'26.76"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +2.19%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  -0.06%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'7.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -5.75%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'0.161"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -6.53%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  -4.28%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'10.29"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.26%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'34.90"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +3.24%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'0.0465"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +3.11%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  -2.37%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'49.83"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -1.43%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +0.06%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'3.31"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +0.44%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.286"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -2.18%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'129.90"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +1.04%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'1.81"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -2.08%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.114"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -1.51%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'16.21"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -4.16%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'3.70"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -2.00%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'2.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -3.44%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'21.41"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -0.54%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'2.48"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +2.74%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'2.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.05%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'2.041.15"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.72%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0318"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +0.14%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.884"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +12.50%  "
$ws.Range("E51").ClearFormats()
